$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Page 1")

# --- Header row 3: relabel/reorder a few columns ---
$ws.Range("J3").Value = "EX"
$ws.Range("L3").Value = "Cuest."
$ws.Range("O3").Value = "Extras"

# --- Column width adjustments ---
$ws.Columns.Item(1).ColumnWidth = 4.42578125
$ws.Columns.Item(4).ColumnWidth = 6.7109375
$ws.Columns.Item(5).ColumnWidth = 6.28515625
$ws.Columns.Item(6).ColumnWidth = 6.42578125
$ws.Columns.Item(7).ColumnWidth = 6.28515625
$ws.Columns.Item(8).ColumnWidth = 5.85546875
$ws.Columns.Item(9).ColumnWidth = 4.85546875
$ws.Columns.Item(10).ColumnWidth = 5.5703125
$ws.Columns.Item(11).ColumnWidth = 5.28515625
$ws.Columns.Item(12).ColumnWidth = 7.42578125

# --- Selection moves to G4 ---
$ws.Range("G4").Select()

# --- Zero out the O column (extra points) for rows 5-122 ---
for ($r = 5; $r -le 122; $r++) {
    $ws.Cells.Item($r, 15).Value = 0
}
